$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    # Preserve the cell's existing style/number-format (e.g. borders, General
    # format) while forcing a literal text assignment, so Excel does not
    # auto-convert percent-looking strings ("72%") into numeric % values.
    $target = $ws.Range($cellRef)
    $target.Copy()
    $ws.Range("ZZ100").PasteSpecial(-4122)  # xlPasteFormats
    $target.NumberFormat = "@"
    $target.Value = $newValue
    $ws.Range("ZZ100").Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("ZZ100").Clear()
}

$ws.Range('E2').Value = '2026-02-14 22:18:46'
$ws.Range('G2').Value = '203 cm'
$ws.Range('N2').Value = '-4.0 °C 21:55 TU'
$ws.Range('E3').Value = '2026-02-14 22:18:48'
$ws.Range('N3').Value = '-8.1 °C 21:59 TU'
$ws.Range('E4').Value = '2026-02-14 22:18:51'
$ws.Range('J4').Value = '998.1 hPa'
$ws.Range('O4').Value = '10.6 °C'
$ws.Range('E5').Value = '2026-02-14 22:18:54'
$ws.Range('I5').Value = '22.7 mm'
$ws.Range('N5').Value = '-7.7 °C 21:57 TU'
$ws.Range('E6').Value = '2026-02-14 22:18:57'
Set-TextValue 'H6' '72%'
$ws.Range('J6').Value = '998.1 hPa'
$ws.Range('E7').Value = '2026-02-14 22:18:59'
Set-TextValue 'H7' '50%'
$ws.Range('J7').Value = '998.3 hPa'
$ws.Range('E8').Value = '2026-02-14 22:19:02'
Set-TextValue 'H8' '60%'
$ws.Range('J8').Value = '998.1 hPa'
$ws.Range('E9').Value = '2026-02-14 22:19:05'
$ws.Range('N9').Value = '8.5 °C 21:54 TU'
$ws.Range('O9').Value = '11.7 °C'
$ws.Range('E10').Value = '2026-02-14 22:19:08'
Set-TextValue 'H10' '74%'
$ws.Range('E11').Value = '2026-02-14 22:19:10'
$ws.Range('E12').Value = '2026-02-14 22:19:13'
$ws.Range('N12').Value = '9.2 °C 21:57 TU'
$ws.Range('O12').Value = '12.0 °C'
$ws.Range('E13').Value = '2026-02-14 22:19:15'
$ws.Range('J13').Value = '1000.9 hPa'
$ws.Range('E14').Value = '2026-02-14 22:19:18'
$ws.Range('E15').Value = '2026-02-14 22:19:21'
$ws.Range('N15').Value = '8.2 °C 21:55 TU'
$ws.Range('E16').Value = '2026-02-14 22:19:24'
Set-TextValue 'H16' '75%'
$ws.Range('K16').Value = '8.8 MJ/m2'
$ws.Range('E17').Value = '2026-02-14 22:19:27'
Set-TextValue 'H17' '66%'
$ws.Range('L17').Value = '69.5 km/h - 58º 21:43 TU'
$ws.Range('N17').Value = '-1.3 °C 21:47 TU'
$ws.Range('O17').Value = '1.6 °C'
$ws.Range('E18').Value = '2026-02-14 22:19:29'
Set-TextValue 'H18' '73%'
$ws.Range('J18').Value = '998.3 hPa'
$ws.Range('E19').Value = '2026-02-14 22:19:32'
Set-TextValue 'H19' '74%'
$ws.Range('L19').Value = '39.2 km/h - 350º 21:40 TU'
$ws.Range('E20').Value = '2026-02-14 22:19:35'
$ws.Range('I20').Value = '5.0 mm'
$ws.Range('N20').Value = '-8.7 °C 21:50 TU'
$ws.Range('E21').Value = '2026-02-14 22:19:37'
Set-TextValue 'H21' '67%'
$ws.Range('J21').Value = '1000.7 hPa'
$ws.Range('E22').Value = '2026-02-14 22:19:40'
Set-TextValue 'H22' '85%'
$ws.Range('O22').Value = '-7.0 °C'
$ws.Range('E23').Value = '2026-02-14 22:19:43'
$ws.Range('I23').Value = '40.7 mm'
$ws.Range('N23').Value = '-8.8 °C 21:56 TU'
$ws.Range('O23').Value = '-6.2 °C'
$ws.Range('E24').Value = '2026-02-14 22:19:45'
$ws.Range('J24').Value = '1002.4 hPa'
$ws.Range('O24').Value = '9.3 °C'
$ws.Range('E25').Value = '2026-02-14 22:19:48'
$ws.Range('I25').Value = '19.4 mm'
$ws.Range('N25').Value = '-7.8 °C 21:36 TU'
$ws.Range('O25').Value = '-4.9 °C'
$ws.Range('E26').Value = '2026-02-14 22:19:51'
$ws.Range('E27').Value = '2026-02-14 22:19:54'
$ws.Range('N27').Value = '-6.2 °C 21:50 TU'
$ws.Range('O27').Value = '-3.3 °C'
$ws.Range('E28').Value = '2026-02-14 22:19:57'
Set-TextValue 'H28' '64%'
$ws.Range('J28').Value = '998.0 hPa'
$ws.Range('L28').Value = '64.1 km/h - 332º 21:56 TU'
$ws.Range('E29').Value = '2026-02-14 22:19:59'
$ws.Range('E30').Value = '2026-02-14 22:20:02'
$ws.Range('J30').Value = '998.0 hPa'
$ws.Range('E31').Value = '2026-02-14 22:20:05'
$ws.Range('J31').Value = '997.1 hPa'
$ws.Range('N31').Value = '7.1 °C 21:59 TU'
$ws.Range('E32').Value = '2026-02-14 22:20:08'
$ws.Range('N32').Value = '1.8 °C 21:59 TU'
$ws.Range('E33').Value = '2026-02-14 22:20:10'
$ws.Range('J33').Value = '1000.4 hPa'
$ws.Range('O33').Value = '3.8 °C'
$ws.Range('E34').Value = '2026-02-14 22:20:13'
$ws.Range('N34').Value = '-5.3 °C 21:54 TU'
$ws.Range('O34').Value = '-2.5 °C'
$ws.Range('E35').Value = '2026-02-14 22:20:16'
$ws.Range('J35').Value = '1004.9 hPa'
$ws.Range('N35').Value = '1.3 °C 21:57 TU'
$ws.Range('E36').Value = '2026-02-14 22:20:18'
$ws.Range('J36').Value = '998.8 hPa'
$ws.Range('L36').Value = '96.8 km/h - 337º 21:55 TU'
$ws.Range('N36').Value = '9.6 °C 21:59 TU'
$ws.Range('E37').Value = '2026-02-14 22:20:21'
$ws.Range('J37').Value = '999.0 hPa'
$ws.Range('E38').Value = '2026-02-14 22:20:23'
$ws.Range('N38').Value = '6.3 °C 21:59 TU'
$ws.Range('O38').Value = '10.0 °C'
$ws.Range('E39').Value = '2026-02-14 22:20:26'
$ws.Range('I39').Value = '14.5 mm'
$ws.Range('N39').Value = '-8.8 °C 21:49 TU'
$ws.Range('E40').Value = '2026-02-14 22:20:29'
Set-TextValue 'H40' '64%'
$ws.Range('J40').Value = '1001.5 hPa'
$ws.Range('E41').Value = '2026-02-14 22:20:32'
$ws.Range('J41').Value = '1000.1 hPa'
$ws.Range('O41').Value = '13.3 °C'
$ws.Range('E42').Value = '2026-02-14 22:20:34'
$ws.Range('E43').Value = '2026-02-14 22:20:37'
$ws.Range('E44').Value = '2026-02-14 22:20:40'
$ws.Range('I44').Value = '37.7 mm'
$ws.Range('N44').Value = '-8.2 °C 21:32 TU'
$ws.Range('E45').Value = '2026-02-14 22:20:42'
$ws.Range('J45').Value = '1007.5 hPa'
$ws.Range('N45').Value = '0.2 °C 21:58 TU'
$ws.Range('O45').Value = '2.8 °C'
$ws.Range('E46').Value = '2026-02-14 22:20:45'
